$wb = $excel.ActiveWorkbook

# Update the source value on the "updates_from_sub" sheet (A_orb row).
# All dependent formulas (Calculations Rough Input!C3/D3/C6/C7/C8/C9,
# EPS!D2/G2/H2/I2/J2/K2/D3/G3/H3/I3/J3/K3) recalculate automatically.
$ws = $wb.Worksheets.Item("updates_from_sub")
$ws.Range("B2").Value = 7.8826646006315837

$excel.CalculateFullRebuild()
